$p = $ppt.ActivePresentation

# 1. Update the table style id on slide 6 (the "SOURCES OF FINANCE" slide's table)
$slide = $p.Slides.Item(6)
$tableShape = $slide.Shapes.Item(2)
$tableShape.Table.ApplyStyle("{87D2A48C-A8AF-4CC3-9805-7B6B8F9CB29D}")
